# Weekly fruit/vegetable price update: a new observation was recorded
# (week of 2022-06-02, serial 44714) and inserted as the new row 59,
# pushing the previously-last three rows down by one (59->60, 60->61,
# 61->62). Row 58 and everything above stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59, shifting rows 59:61 down to 60:62 and
# carrying the date-format style (s=2) of column D down with them.
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new weekly data point.
$ws.Cells.Item(59, 1).Value  = 4
$ws.Cells.Item(59, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value  = "Los Lagos"
$ws.Cells.Item(59, 4).Value  = 44714
$ws.Cells.Item(59, 5).Value  = 10
$ws.Cells.Item(59, 6).Value  = 100112043
$ws.Cells.Item(59, 7).Value  = "Pepino dulce"
$ws.Cells.Item(59, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(59, 9).Value  = "Especial"
$ws.Cells.Item(59, 10).Value = 30
$ws.Cells.Item(59, 11).Value = 21000
$ws.Cells.Item(59, 12).Value = 21000
$ws.Cells.Item(59, 13).Value = 21000
$ws.Cells.Item(59, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value = 1167
$ws.Cells.Item(59, 17).Value = 18
$ws.Cells.Item(59, 18).Value = "Hortaliza"
